# Weekly data update: insert a new weekly observation as row 29, pushing the
# existing rows 29-39 down to 30-40 (dimension grows from A1:R39 to A1:R40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 29. This shifts rows 29-39
# down to rows 30-40 (and copies the row-29 formatting, including the date
# style on column D, into the new row).
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new weekly record. The descriptive
# columns (market id/name/region, category id/name, variety, quality,
# unit, origin, kg-or-units, classification) are identical to the
# surrounding rows, so copy them from the row just below (old row 29,
# now row 30). Only the date and the price-related columns differ.
$ws.Cells.Item(29, 1).Value  = $ws.Cells.Item(30, 1).Value()   # Mercado ID
$ws.Cells.Item(29, 2).Value  = $ws.Cells.Item(30, 2).Value()   # Mercado
$ws.Cells.Item(29, 3).Value  = $ws.Cells.Item(30, 3).Value()   # Región
$ws.Cells.Item(29, 4).Value  = 44524                           # Fecha
$ws.Cells.Item(29, 5).Value  = $ws.Cells.Item(30, 5).Value()   # Codreg
$ws.Cells.Item(29, 6).Value  = $ws.Cells.Item(30, 6).Value()   # Categoría ID
$ws.Cells.Item(29, 7).Value  = $ws.Cells.Item(30, 7).Value()   # Categoría
$ws.Cells.Item(29, 8).Value  = $ws.Cells.Item(30, 8).Value()   # Variedad
$ws.Cells.Item(29, 9).Value  = $ws.Cells.Item(30, 9).Value()   # Calidad
$ws.Cells.Item(29, 10).Value = 16                               # Volumen
$ws.Cells.Item(29, 11).Value = 9000                             # Precio mínimo
$ws.Cells.Item(29, 12).Value = 10000                            # Precio máximo
$ws.Cells.Item(29, 13).Value = 9500                             # Precio promedio ponderado
$ws.Cells.Item(29, 14).Value = $ws.Cells.Item(30, 14).Value()  # Unidad de comercialización
$ws.Cells.Item(29, 15).Value = $ws.Cells.Item(30, 15).Value()  # Origen
$ws.Cells.Item(29, 16).Value = 3167                             # Precio $/Kg
$ws.Cells.Item(29, 17).Value = 3                                # Kg o Unidades
$ws.Cells.Item(29, 18).Value = $ws.Cells.Item(30, 18).Value()  # Clasificación
